# The deck ships two theme parts:
#   - theme1.xml  -> bound to the (single) Slide Master, i.e. the theme
#                    that actually drives the look of the deck ("Integral").
#   - theme2.xml  -> bound to the Notes Master ("Office Theme", the plain
#                    PowerPoint default palette).
#
# The authored edit swaps the two themes' content: the Slide Master's
# theme becomes the plain "Office Theme" colour palette (what used to
# live in theme2.xml) while the Notes Master conceptually keeps the
# "Integral" colours. The two themes already share an identical font
# scheme and format scheme (fills/lines/effects) - only the 12 theme
# colours (and the cosmetic theme/colour-scheme names) differ between
# them, so the visible, reproducible part of this edit is: re-point the
# Slide Master's theme colours at the "Office Theme" palette.

$p = $ppt.ActivePresentation

# Helper matching VBA's RGB() - packs R,G,B into the BGR-ordered long
# that PowerPoint's ColorFormat/ThemeColor .RGB property expects.
function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# The theme colour scheme is reachable from any slide (it is shared by
# the whole deck via the Slide Master) - use the first slide as anchor.
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

# Standard Office theme palette, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeTheme = @(
    (RGB 0x00 0x00 0x00),  # dk1      000000
    (RGB 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (RGB 0x44 0x54 0x6A),  # dk2      44546A
    (RGB 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (RGB 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (RGB 0xED 0x7D 0x31),  # accent2  ED7D31
    (RGB 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (RGB 0xFF 0xC0 0x00),  # accent4  FFC000
    (RGB 0x44 0x72 0xC4),  # accent5  4472C4
    (RGB 0x70 0xAD 0x47),  # accent6  70AD47
    (RGB 0x05 0x63 0xC1),  # hlink    0563C1
    (RGB 0x95 0x4F 0x72)   # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeTheme[$i - 1]
}
